$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Target cluster: ECs)
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 4.619088000000001
$ws.Range("H2").Value = 13.857264
$ws.Range("K2").Value = 2
$ws.Range("M2").Value = 8.915438
$ws.Range("N2").Value = 17.830876
$ws.Range("O2").Value = 0.1330966619879936
$ws.Range("P2").Value = 0.1044454436141807
$ws.Range("Q2").Value = 41.18119268054401
$ws.Range("R2").Value = 247.087156083264
$ws.Range("S2").Value = 0.1330966619879936
$ws.Range("T2").Value = 0.1044454436141807

# Row 3 (Target cluster: FAPs)
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 4.619088000000001
$ws.Range("H3").Value = 13.857264
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 28.70072166666667
$ws.Range("N3").Value = 86.10216500000001
$ws.Range("O3").Value = 0.42846691889729
$ws.Range("P3").Value = 0.5043486825642434
$ws.Range("Q3").Value = 132.57115904184
$ws.Range("R3").Value = 1193.14043137656
$ws.Range("S3").Value = 0.42846691889729
$ws.Range("T3").Value = 0.5043486825642434

# Row 4 (Target cluster: M1)
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 4.619088000000001
$ws.Range("H4").Value = 13.857264
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.2432513333333333
$ws.Range("N4").Value = 0.729754
$ws.Range("O4").Value = 0.003631446990130538
$ws.Range("P4").Value = 0.004274578560202137
$ws.Range("Q4").Value = 1.123599314784
$ws.Range("R4").Value = 10.112393833056
$ws.Range("S4").Value = 0.003631446990130538
$ws.Range("T4").Value = 0.004274578560202137

# Row 5 (Target cluster: M2)
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 4.619088000000001
$ws.Range("H5").Value = 13.857264
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.2705943333333333
$ws.Range("N5").Value = 0.8117829999999999
$ws.Range("O5").Value = 0.004039644773429317
$ws.Range("P5").Value = 0.004755068430370469
$ws.Range("Q5").Value = 1.249899037968
$ws.Range("R5").Value = 11.249091341712
$ws.Range("S5").Value = 0.004039644773429317
$ws.Range("T5").Value = 0.004755068430370469

# Row 6 (Target cluster: Neutro)
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 4.619088000000001
$ws.Range("H6").Value = 13.857264
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 7.535582333333333
$ws.Range("N6").Value = 22.606747
$ws.Range("O6").Value = 0.1124970926501157
$ws.Range("P6").Value = 0.1324203992607289
$ws.Range("Q6").Value = 34.80751792891201
$ws.Range("R6").Value = 313.267661360208
$ws.Range("S6").Value = 0.1124970926501157
$ws.Range("T6").Value = 0.1324203992607289

# Row 7 (Target cluster: sCs)
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 4.619088000000001
$ws.Range("H7").Value = 13.857264
$ws.Range("K7").Value = 2
$ws.Range("M7").Value = 21.3190975
$ws.Range("N7").Value = 42.638195
$ws.Range("O7").Value = 0.3182682347010409
$ws.Range("P7").Value = 0.2497558275702743
$ws.Range("Q7").Value = 98.47478743308
$ws.Range("R7").Value = 590.84872459848
$ws.Range("S7").Value = 0.3182682347010409
$ws.Range("T7").Value = 0.2497558275702743
